$d = $word.ActiveDocument

# --- 1. Remove the "_GoBack" bookmark from the paragraph that currently
#        holds it (the paragraph becomes an ordinary empty paragraph). ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Append new content at the end of the document body. ---
$tail = $d.Paragraphs.Last.Range
$tail.Collapse(0)

# Helper: start a brand-new empty paragraph at $tail and land the range
# cursor *inside* it (past the freshly inserted paragraph mark). New
# paragraphs inherit the style of whatever preceded them, so pin it back
# to "Normal" by default (callers can override afterwards).
function New-TailParagraph($r) {
    $r.InsertParagraphAfter()
    $r.Move(1, 1) | Out-Null
    $d.Paragraphs.Last.Style = "Normal"
}

# Seven blank paragraphs.
for ($i = 0; $i -lt 7; $i++) {
    New-TailParagraph $tail
}

# Heading: "Séance du 07/01/19"
New-TailParagraph $tail
$tail.InsertAfter("Séance du 07/01/19")
$d.Paragraphs.Last.Style = "Heading 2"
$tail.Collapse(0)

# Paragraph: movement / Bluetooth piloting text + video intro
New-TailParagraph $tail
$tail.InsertAfter("Maintenant que nous avons réglé le système de tir, nous nous occupons du déplacement du robot. A partir des anciens TD (TD 8), nous avons pu assez rapidement réussir à avoir un bon déplacement (avant, arrière, gauche, droite), malgré la limite de la portée du Bluetooth. Nous pensons tout de même que cela ne posera pas de problème avec les dimensions du terrain. Voici une vidéo de pilotage de notre robot :")
$tail.Collapse(0)

# Paragraph: video link
New-TailParagraph $tail
$tail.InsertAfter("https://www.youtube.com/watch?v=n10RO3F-3z0")
$tail.Collapse(0)

# Paragraph: program changes / module at the front
New-TailParagraph $tail
$tail.InsertAfter("Nous avons donc apporté des modifications aux programmes pour avoir un robot adapté à nos besoins, en inversant par exemple l’avant et l’arrière du robot. En effet, nous réfléchissons donc à la disposition des éléments pour que notre robot soit efficace, malgré le fait que nous ne puissions pas changer la forme étant donné que nous utilisons un robot de TD. Nous allons donc créer un module à l’avant afin de contrôler et tirer la balle, car le châssis est trop bas pour que la balle passe dessous.")
$tail.Collapse(0)

# Paragraph: servomotors / Bluetooth module issue
New-TailParagraph $tail
$tail.InsertAfter("Nous continuons également de chercher des solutions pour utiliser plusieurs servomoteurs, car en testant avec un programme simple nous n’arrivons pas à en faire marcher deux sur la même carte, alors que nous avons besoin qu’ils soient sur le même module Bluetooth que les moteurs.")
$tail.Collapse(0)

# Paragraph: battery power solution (ends with the relocated bookmark)
New-TailParagraph $tail
$tail.InsertAfter("Nous allons donc essayer une solution qui consiste à les alimenter directement depuis les piles pour avoir un courant supérieur ce qui permettrai de faire fonctionner les deux moteurs. Nous allons réaliser les branchements dans la semaine, et si cela fonctionne nous mettrons en place un module à l’avant du robot adapté aux deux servomoteurs que nous avons.")
$tail.Collapse(0)

$d.Bookmarks.Add("_GoBack", $tail)

# Final trailing blank paragraph.
New-TailParagraph $tail
